$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_ALC.Range("H70").Value = 937.05884
$ws_ALC.Range("I70").Value = 945.2
$ws_ALC.Range("J70").Value = 914.44446
$ws_ALC.Range("K70").Value = 2835.6
$ws_ALC.Range("L70").Value = 2743.33338
$ws_ALC.Range("M70").Value = -2565.6
$ws_ALC.Range("N70").Value = -3283.33338
$ws_ALC.Range("H73").Value = 937.05884
$ws_ALC.Range("I73").Value = 945.2
$ws_ALC.Range("J73").Value = 914.44446
$ws_ALC.Range("K73").Value = 2835.6
$ws_ALC.Range("L73").Value = 2743.33338
$ws_ALC.Range("M73").Value = -1899.6
$ws_ALC.Range("N73").Value = -4615.33338
$ws_ALC.Range("H121").Value = 3239
$ws_ALC.Range("J121").Value = 6605
$ws_ALC.Range("L121").Value = 19815
$ws_ALC.Range("N121").Value = -23309
$ws_ALC.Range("H137").Value = 2482.35
$ws_ALC.Range("I137").Value = 4228.5557
$ws_ALC.Range("J137").Value = 1053.6364
$ws_ALC.Range("K137").Value = 12685.6671
$ws_ALC.Range("L137").Value = 3160.9092
$ws_ALC.Range("M137").Value = -10135.6671
$ws_ALC.Range("N137").Value = -8260.9092
$ws_ALC.Range("H141").Value = 9718.200000000001
$ws_ALC.Range("I141").Value = 3161.4285
$ws_ALC.Range("J141").Value = 25017.334
$ws_ALC.Range("K141").Value = 9484.2855
$ws_ALC.Range("L141").Value = 75052.00199999999
$ws_ALC.Range("M141").Value = -4304.2855
$ws_ALC.Range("N141").Value = -85412.00199999999
$ws_ARM.Range("H32").Value = 5283.43
$ws_ARM.Range("I32").Value = 4561.5605
$ws_ARM.Range("J32").Value = 12582.333
$ws_ARM.Range("K32").Value = 4561.5605
$ws_ARM.Range("L32").Value = 12582.333
$ws_ARM.Range("M32").Value = -4274.5605
$ws_ARM.Range("N32").Value = -13156.333
$ws_ARM.Range("H45").Value = 1431.375
$ws_ARM.Range("I45").Value = 1233.7391
$ws_ARM.Range("J45").Value = 1936.4445
$ws_ARM.Range("K45").Value = 1233.7391
$ws_ARM.Range("L45").Value = 1936.4445
$ws_ARM.Range("M45").Value = -856.7391
$ws_ARM.Range("N45").Value = -2690.4445
$ws_ARM.Range("H61").Value = 4786.8096
$ws_ARM.Range("I61").Value = 6501.846
$ws_ARM.Range("J61").Value = 1999.875
$ws_ARM.Range("K61").Value = 6501.846
$ws_ARM.Range("L61").Value = 1999.875
$ws_ARM.Range("M61").Value = -6289.846
$ws_ARM.Range("N61").Value = -2423.875
$ws_ARM.Range("H132").Value = 976700.6
$ws_ARM.Range("I132").Value = 2509762.8
$ws_ARM.Range("J132").Value = 5761.2
$ws_ARM.Range("K132").Value = 7529288.399999999
$ws_ARM.Range("L132").Value = 17283.6
$ws_ARM.Range("M132").Value = -7526758.399999999
$ws_ARM.Range("N132").Value = -22343.6
$ws_ARM.Range("H136").Value = 4786.8096
$ws_ARM.Range("I136").Value = 6501.846
$ws_ARM.Range("J136").Value = 1999.875
$ws_ARM.Range("K136").Value = 19505.538
$ws_ARM.Range("L136").Value = 5999.625
$ws_ARM.Range("M136").Value = -16955.538
$ws_ARM.Range("N136").Value = -11099.625
$ws_BSM.Range("H62").Value = 40000
$ws_BSM.Range("J62").Value = 40000
$ws_BSM.Range("L62").Value = 40000
$ws_BSM.Range("N62").Value = -41372
$ws_BSM.Range("H65").Value = 40000
$ws_BSM.Range("J65").Value = 40000
$ws_BSM.Range("L65").Value = 120000
$ws_BSM.Range("N65").Value = -126864
$ws_BSM.Range("H92").Value = 26500
$ws_BSM.Range("J92").Value = 26500
$ws_BSM.Range("L92").Value = 26500
$ws_BSM.Range("N92").Value = -31492
$ws_BSM.Range("H100").Value = 20000
$ws_BSM.Range("J100").Value = 20000
$ws_BSM.Range("L100").Value = 20000
$ws_BSM.Range("N100").Value = -22164
$ws_BSM.Range("H139").Value = 64970
$ws_BSM.Range("J139").Value = 64970
$ws_BSM.Range("L139").Value = 64970
$ws_BSM.Range("N139").Value = -75250
$ws_CRP.Range("H4").Value = 9397.223
$ws_CRP.Range("J4").Value = 9821.875
$ws_CRP.Range("L4").Value = 9821.875
$ws_CRP.Range("N4").Value = -10045.875
$ws_CRP.Range("H58").Value = 5816.4644
$ws_CRP.Range("I58").Value = 3188.6191
$ws_CRP.Range("K58").Value = 3188.6191
$ws_CRP.Range("M58").Value = -2985.6191
$ws_CRP.Range("H74").Value = 21200.25
$ws_CRP.Range("J74").Value = 21200.25
$ws_CRP.Range("L74").Value = 21200.25
$ws_CRP.Range("N74").Value = -22948.25
$ws_CRP.Range("H77").Value = 21200.25
$ws_CRP.Range("J77").Value = 21200.25
$ws_CRP.Range("L77").Value = 63600.75
$ws_CRP.Range("N77").Value = -72336.75
$ws_CRP.Range("H134").Value = 3416
$ws_CRP.Range("J134").Value = 3970
$ws_CRP.Range("L134").Value = 11910
$ws_CRP.Range("N134").Value = -16980
$ws_CRP.Range("H136").Value = 5816.4644
$ws_CRP.Range("I136").Value = 3188.6191
$ws_CRP.Range("K136").Value = 9565.8573
$ws_CRP.Range("M136").Value = -7015.8573
$ws_CUL.Range("H4").Value = 5843.1113
$ws_CUL.Range("I4").Value = 8514.666999999999
$ws_CUL.Range("K4").Value = 25544.001
$ws_CUL.Range("M4").Value = -25432.001
$ws_CUL.Range("H5").Value = 1120.8334
$ws_CUL.Range("I5").Value = 832
$ws_CUL.Range("J5").Value = 2565
$ws_CUL.Range("K5").Value = 2496
$ws_CUL.Range("L5").Value = 7695
$ws_CUL.Range("M5").Value = -2384
$ws_CUL.Range("N5").Value = -7919
$ws_CUL.Range("H131").Value = 1035.9767
$ws_CUL.Range("I131").Value = 616.6667
$ws_CUL.Range("J131").Value = 1103.973
$ws_CUL.Range("K131").Value = 1850.0001
$ws_CUL.Range("L131").Value = 3311.919
$ws_CUL.Range("M131").Value = 3189.9999
$ws_CUL.Range("N131").Value = -13391.919
$ws_CUL.Range("H135").Value = 1120.8334
$ws_CUL.Range("I135").Value = 832
$ws_CUL.Range("J135").Value = 2565
$ws_CUL.Range("K135").Value = 7488
$ws_CUL.Range("L135").Value = 23085
$ws_CUL.Range("M135").Value = -4953
$ws_CUL.Range("N135").Value = -28155
$ws_GSM.Range("H5").Value = 5
$ws_GSM.Range("J5").Value = 5
$ws_GSM.Range("L5").Value = 5
$ws_GSM.Range("N5").Value = -229
$ws_GSM.Range("H98").Value = 8577.200000000001
$ws_GSM.Range("J98").Value = 8577.200000000001
$ws_GSM.Range("L98").Value = 8577.200000000001
$ws_GSM.Range("N98").Value = -14567.2
$ws_GSM.Range("H132").Value = 2981563.8
$ws_GSM.Range("I132").Value = 10421872
$ws_GSM.Range("J132").Value = 5440.6
$ws_GSM.Range("K132").Value = 31265616
$ws_GSM.Range("L132").Value = 16321.8
$ws_GSM.Range("M132").Value = -31263086
$ws_GSM.Range("N132").Value = -21381.8
$ws_LTW.Range("H2").Value = 2148068.5
$ws_LTW.Range("J2").Value = 2148068.5
$ws_LTW.Range("L2").Value = 2148068.5
$ws_LTW.Range("N2").Value = -2148292.5
$ws_LTW.Range("H136").Value = 2712.8572
$ws_LTW.Range("I136").Value = 2000
$ws_LTW.Range("J136").Value = 2831.6667
$ws_LTW.Range("K136").Value = 6000
$ws_LTW.Range("L136").Value = 8495.000100000001
$ws_LTW.Range("M136").Value = -3450
$ws_LTW.Range("N136").Value = -13595.0001
$ws_WVR.Range("H2").Value = 148593860
$ws_WVR.Range("I2").Value = 20005500
$ws_WVR.Range("J2").Value = 200029200
$ws_WVR.Range("K2").Value = 20005500
$ws_WVR.Range("L2").Value = 200029200
$ws_WVR.Range("M2").Value = -20005388
$ws_WVR.Range("N2").Value = -200029424
$ws_WVR.Range("H81").Value = 1925.579
$ws_WVR.Range("I81").Value = 2038.6923
$ws_WVR.Range("J81").Value = 1680.5
$ws_WVR.Range("K81").Value = 4077.3846
$ws_WVR.Range("L81").Value = 3361
$ws_WVR.Range("M81").Value = -3016.3846
$ws_WVR.Range("N81").Value = -5483
$ws_WVR.Range("H84").Value = 1925.579
$ws_WVR.Range("I84").Value = 2038.6923
$ws_WVR.Range("J84").Value = 1680.5
$ws_WVR.Range("K84").Value = 20386.923
$ws_WVR.Range("L84").Value = 16805
$ws_WVR.Range("M84").Value = -15082.923
$ws_WVR.Range("N84").Value = -27413

Write-Output "Applied 183 cell updates"
